$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the "To do" column with 5 new Kanban cards (rows 36-40) ---
# Clone the formatting of the last existing row (35) down across the new rows
# first, so the new cells inherit the same look (left/right/top thin borders
# on column C, top-thin-only on A/B) without creating throwaway style
# records.
$ws.Range("A35:C35").Copy()
$ws.Range("A36:C40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column C on these new rows has no right-hand neighbour cell needing a
# shared edge, so drop the right border that was copied along with the rest
# of the "left+right+top" look, leaving "left+top" only.
$ws.Range("C36:C40").Borders.Item(10).LineStyle = -4142  # xlEdgeRight, xlLineStyleNone

# Task text for the new cards.
$ws.Range("B36").Value = "Crear el componente Catalogo.js"
$ws.Range("B37").Value = "Configurar las rutas en React Router"
$ws.Range("B38").Value = "Agregar un enlace al catálogo en la barra de navegación"
$ws.Range("B39").Value = "Agregar estilos con Bootstrap"
$ws.Range("B40").Value = "Probar la página del catálogo en el navegador"

# "Catalogo.js" inside the B36 text is rendered in a distinct run (smaller,
# fixed-width-ish font) to match the rest of the sheet's convention for
# highlighting file/module names.
$run = $ws.Range("B36").Characters(21, 11)
$run.Font.Size = 10
$run.Font.Name = "Arial Unicode MS"

# --- Grow the table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C40"))

# --- Update the view so it reflects the extended list ---
$ws.Range("B40").Select()
